$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156 (shifts existing rows 156:184 down to 157:185)
$ws.Rows("156:156").Insert()

# Populate the newly inserted row 156 with the new weekly price record
$ws.Cells.Item(156, 1).Value = 7
$ws.Cells.Item(156, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(156, 3).Value = "Ñuble"
$ws.Cells.Item(156, 4).Value = 44504
$ws.Cells.Item(156, 5).Value = 16
$ws.Cells.Item(156, 6).Value = 100112023
$ws.Cells.Item(156, 7).Value = "Brócoli"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 300
$ws.Cells.Item(156, 11).Value = 650
$ws.Cells.Item(156, 12).Value = 700
$ws.Cells.Item(156, 13).Value = 675
$ws.Cells.Item(156, 14).Value = "$/unidad"
$ws.Cells.Item(156, 15).Value = "Región del Maule"
$ws.Cells.Item(156, 16).Value = 675
$ws.Cells.Item(156, 17).Value = 1
$ws.Cells.Item(156, 18).Value = "Hortaliza"
